$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the screened client emails (grading finished per subject) ---

# Row 2 (SN 1) and Row 3 (SN 2) belong to the same recipient: odenadoma@gmail.com
$ws.Range("I2").Value = "odenadoma@gmail.com"
$ws.Range("I3").Value = "odenadoma@gmail.com"

# Row 4 (SN 3) belongs to a different recipient: greyspades99@gmail.com
$ws.Range("I4").Value = "greyspades99@gmail.com"

# I2 already had a mailto: hyperlink - repoint it at the corrected address.
$ws.Range("I2").Hyperlinks.Item(1).Address = "mailto:odenadoma@gmail.com"

# I3 / I4 did not have hyperlinks yet - add them now.
$ws.Range("I3").Hyperlinks.Add($ws.Range("I3"), "mailto:odenadoma@gmail.com") | Out-Null
$ws.Range("I3").Font.Underline = 2

$ws.Range("I4").Hyperlinks.Add($ws.Range("I4"), "mailto:greyspades99@gmail.com") | Out-Null
$ws.Range("I4").Font.Underline = 2

# Move the live selection to where the grader left off.
$ws.Range("K18").Select() | Out-Null
